$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 37

# Columns A (date-looking) and C (numeric-looking) must be forced to Text
# format *before* assigning their values, otherwise Excel auto-detects them
# as a date serial / a number instead of keeping the literal text that the
# source data feed writes (t="str" cells in the OOXML).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025-10-23"

$ws.Cells.Item($newRow, 2).Value = "Pick 3"

$ws.Cells.Item($newRow, 3).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).Value = "251023"

$ws.Cells.Item($newRow, 4).Value = "5-2-7"

$ws.Cells.Item($newRow, 5).Value = "2025-10-23T21:38:22.290+04:00"
